# Team_Project Diary.xlsx - work log update for Week 5 / Week 6
# "Updated work log with hours for week 5/6"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Week 5 / Week 6 sheets: add the new "data extractor" work items.
# The order of the writes below matters: it reproduces the order the
# shared-string table entries were originally appended in (indices 74-77).
# ---------------------------------------------------------------------------
$wk5 = $wb.Worksheets.Item("Week 5")
$wk6 = $wb.Worksheets.Item("Week 6")

# 1) Week 5, B18 -> new shared string #74
$wk5.Range("B18").Value = "Huge amount of work improving the data extractor program"

# 2) Week 6, B6 -> new shared string #75 (replaces the old "Week 5 (26th..." label)
$wk6.Range("B6").Value = "Week 5 (31th Oct – 6th November)"

# 3) Week 6, B18 -> new shared string #76
$wk6.Range("B18").Value = "Compiling output of ""data_extractor.py"""

# 4) Week 5, B19 -> new shared string #77
$wk5.Range("B19").Value = "Compiling output of ""data_extractor.py"" for high/medium/random groups (low group failed to work)"
$wk5.Range("C19").Value = 4

# 5) Week 5, B21 -> reuses the existing "Friday meeting" shared string
$wk5.Range("B21").Value = "Friday meeting"
$wk5.Range("C21").Value = 1

# ---------------------------------------------------------------------------
# Restore the selections left behind by the author on each sheet
# ---------------------------------------------------------------------------
$wk5.Range("B25").Select()
$wk6.Range("B18").Select()

# Week 5 is the sheet that was active/selected when the workbook was saved
$wk5.Activate()
